$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns for rows 2-51
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.356.25"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.042.41"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.14"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.649"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.78"
$ws.Range("E7").Value = "  +9.65%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.400"
$ws.Range("E9").Value = "  +10.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.35"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  +10.19%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.918"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.74"
$ws.Range("E14").Value = "  +25.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.78"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.343.09"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.68"
$ws.Range("E17").Value = "  +6.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.045.00"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.246.21"
$ws.Range("E19").Value = "  +4.11%  "
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +7.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.49"
$ws.Range("E22").Value = "  +5.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.70"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  +4.88%  "
$ws.Range("E27").Value = "  +4.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.06"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.10"
$ws.Range("E29").Value = "  +3.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.128"
$ws.Range("E30").Value = "  +29.92%  "
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.18"
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("E33").Value = "  +4.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0629"
$ws.Range("E34").Value = "  +4.73%  "
$ws.Range("E35").Value = "  +5.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.38"
$ws.Range("E36").Value = "  +12.04%  "
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("E40").Value = "  +30.77%  "
$ws.Range("E41").Value = "  +6.24%  "
$ws.Range("E42").Value = "  +7.76%  "
$ws.Range("E43").Value = "  +5.40%  "
$ws.Range("E44").Value = "  +4.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.37"
$ws.Range("E45").Value = "  +4.85%  "
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "95.44"
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.81"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.390.18"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.94"
$ws.Range("E51").Value = "  +0.65%  "
